# Insert two new rows after row 480 (at position 481), shifting existing
# rows 481..583 down to 483..585, and populate the two new rows with the
# new data points described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(481).Insert()
$ws.Rows.Item(481).Insert()

# New row 481
$ws.Cells.Item(481,1).Value = 9
$ws.Cells.Item(481,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(481,3).Value = "Metropolitana"
$ws.Cells.Item(481,4).Value = 44798
$ws.Cells.Item(481,5).Value = 13
$ws.Cells.Item(481,6).Value = 100112031
$ws.Cells.Item(481,7).Value = "Poroto verde"
$ws.Cells.Item(481,8).Value = "Magnum"
$ws.Cells.Item(481,9).Value = "Primera"
$ws.Cells.Item(481,10).Value = 65
$ws.Cells.Item(481,11).Value = 32000
$ws.Cells.Item(481,12).Value = 33000
$ws.Cells.Item(481,13).Value = 32385
$ws.Cells.Item(481,14).Value = "$/malla 25 kilos"
$ws.Cells.Item(481,15).Value = "Perú"
$ws.Cells.Item(481,16).Value = 1295
$ws.Cells.Item(481,17).Value = 25
$ws.Cells.Item(481,18).Value = "Hortaliza"

# New row 482
$ws.Cells.Item(482,1).Value = 9
$ws.Cells.Item(482,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(482,3).Value = "Metropolitana"
$ws.Cells.Item(482,4).Value = 44798
$ws.Cells.Item(482,5).Value = 13
$ws.Cells.Item(482,6).Value = 100112031
$ws.Cells.Item(482,7).Value = "Poroto verde"
$ws.Cells.Item(482,8).Value = "Sin especificar"
$ws.Cells.Item(482,9).Value = "Primera"
$ws.Cells.Item(482,10).Value = 38
$ws.Cells.Item(482,11).Value = 37000
$ws.Cells.Item(482,12).Value = 38000
$ws.Cells.Item(482,13).Value = 37526
$ws.Cells.Item(482,14).Value = "$/malla 25 kilos"
$ws.Cells.Item(482,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(482,16).Value = 1501
$ws.Cells.Item(482,17).Value = 25
$ws.Cells.Item(482,18).Value = "Hortaliza"
